# Calling app now calls its created class: "word type" -> "word_type" header,
# and the unused "plural" column (F) is cleared out (header + the one stray
# "ktar" value in F3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header B1 from "word type" to "word_type"
$ws.Range("B1").Value = "word_type"

# Clear the "plural" column header and its only data value
$ws.Range("F1").ClearContents()
$ws.Range("F3").ClearContents()

# Leave selection on E8 (matches the final cursor position recorded for this edit)
$ws.Range("E8").Select()
